# Update the "Path to Graduation" schedule with the newly uploaded course list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# --- Fall 2022 block (rows 4-9) ---
$ws.Range("A4").Value = "GEOL 1011K"
$ws.Range("B4").Value = 4
$ws.Range("E4").Value = "CPSC 4899"

$ws.Range("A5").Value = "POLS 1101"
$ws.Range("C5").Value = "CPSC 4135"
$ws.Range("D5").Value = 3

$ws.Range("A6").Value = "DSCI 3111"
$ws.Range("C6").Value = "CYBR 4145"

$ws.Range("A7").Value = "CPSC 3121"

$ws.Range("A8").Value = "CPSC 3415"
$ws.Range("B8").Value = 1

$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()

# --- Fall 2023 block (rows 13-15): drop CPSC 4205 ---
$ws.Range("A15").ClearContents()
$ws.Range("B15").ClearContents()
